# Atualização de bases das ligas, do dia: 17-03-2024 às 10:24
# The match records that were stored on a couple of adjacent rows had been
# mixed up. For each pair of rows below, swap the whole record (every
# column from B/HomeTeam-id through AC/PL_AhUnder) between the two rows,
# while leaving column A (the sequential "id") untouched on its row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B (2) .. AC (29)
$firstCol = 2
$lastCol = 29

function Swap-RowData([int]$row1, [int]$row2) {
    for ($col = $firstCol; $col -le $lastCol; $col++) {
        $cell1 = $ws.Cells.Item($row1, $col)
        $cell2 = $ws.Cells.Item($row2, $col)
        $v1 = $cell1.Value2
        $v2 = $cell2.Value2
        # Skip cells that already hold the same value on both rows so we
        # don't needlessly re-serialize a float (avoids introducing binary
        # floating point noise into values that do not actually change).
        if ("$v1" -ne "$v2") {
            $cell1.Value2 = $v2
            $cell2.Value2 = $v1
        }
    }
}

$pairs = @(
    @(8, 9),
    @(13, 14),
    @(18, 19),
    @(26, 27),
    @(30, 31),
    @(47, 48),
    @(54, 55)
)

foreach ($pair in $pairs) {
    Swap-RowData $pair[0] $pair[1]
}
